# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell without letting Excel reinterpret
# numeric-looking strings (e.g. "0.9970", "10.70") as real numbers, which
# would silently drop significant trailing/leading zeros. We briefly mark the
# cell as Text, assign the literal string, then restore the Normal style so the
# cell's formatting stays exactly as it was before the edit.
function Set-TextValue($cell, $text) {
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

Set-TextValue $ws.Range("D2") '26.927.50'
Set-TextValue $ws.Range("E2") '  -1.68%  '
Set-TextValue $ws.Range("D3") '1.811.91'
Set-TextValue $ws.Range("E3") '  -0.78%  '
Set-TextValue $ws.Range("D4") '0.9982'
Set-TextValue $ws.Range("E4") '  -0.33%  '
Set-TextValue $ws.Range("D5") '309.47'
Set-TextValue $ws.Range("E5") '  -1.47%  '
Set-TextValue $ws.Range("D6") '0.9973'
Set-TextValue $ws.Range("E6") '  -0.34%  '
Set-TextValue $ws.Range("D7") '0.4619'
Set-TextValue $ws.Range("E7") '  +3.36%  '
Set-TextValue $ws.Range("D8") '0.3738'
Set-TextValue $ws.Range("E8") '  -0.53%  '
Set-TextValue $ws.Range("D9") '0.07297'
Set-TextValue $ws.Range("E9") '  -2.91%  '
Set-TextValue $ws.Range("D10") '0.8634'
Set-TextValue $ws.Range("E10") '  -3.24%  '
Set-TextValue $ws.Range("E11") '  -3.07%  '
Set-TextValue $ws.Range("D12") '1.739.38'
Set-TextValue $ws.Range("E12") '  -4.87%  '
Set-TextValue $ws.Range("D13") '5.338'
Set-TextValue $ws.Range("E13") '  -1.44%  '
Set-TextValue $ws.Range("D14") '6.515'
Set-TextValue $ws.Range("E14") '  -3.59%  '
Set-TextValue $ws.Range("D15") '0.07029'
Set-TextValue $ws.Range("E15") '  -1.14%  '
Set-TextValue $ws.Range("D16") '91.44'
Set-TextValue $ws.Range("E16") '  -2.91%  '
Set-TextValue $ws.Range("D17") '0.9970'
Set-TextValue $ws.Range("E17") '  -0.49%  '
Set-TextValue $ws.Range("D18") '0.000008663'
Set-TextValue $ws.Range("E18") '  -1.66%  '
Set-TextValue $ws.Range("D19") '0.9974'
Set-TextValue $ws.Range("E19") '  -0.30%  '
Set-TextValue $ws.Range("D20") '14.73'
Set-TextValue $ws.Range("E20") '  -3.20%  '
Set-TextValue $ws.Range("D21") '26.948.46'
Set-TextValue $ws.Range("E21") '  -1.65%  '
Set-TextValue $ws.Range("D22") '5.297'
Set-TextValue $ws.Range("E22") '  +0.62%  '
Set-TextValue $ws.Range("D23") '10.70'
Set-TextValue $ws.Range("E23") '  -2.08%  '
Set-TextValue $ws.Range("D24") '1.963.92'
Set-TextValue $ws.Range("E24") '  -4.48%  '
Set-TextValue $ws.Range("D25") '1.907'
Set-TextValue $ws.Range("E25") '  -3.47%  '
Set-TextValue $ws.Range("D26") '150.73'
Set-TextValue $ws.Range("E26") '  -0.48%  '
Set-TextValue $ws.Range("D27") '18.32'
Set-TextValue $ws.Range("E27") '  -1.51%  '
Set-TextValue $ws.Range("D28") '2.153'
Set-TextValue $ws.Range("E28") '  -9.24%  '
Set-TextValue $ws.Range("E29") '  -2.24%  '
Set-TextValue $ws.Range("D30") '114.71'
Set-TextValue $ws.Range("E30") '  -2.52%  '
Set-TextValue $ws.Range("D31") '0.08885'
Set-TextValue $ws.Range("E31") '  +0.68%  '
Set-TextValue $ws.Range("D32") '0.7659'
Set-TextValue $ws.Range("E32") '  -2.51%  '
Set-TextValue $ws.Range("D33") '1.173'
Set-TextValue $ws.Range("E33") '  -2.46%  '
Set-TextValue $ws.Range("D34") '4.469'
Set-TextValue $ws.Range("E34") '  -1.37%  '
Set-TextValue $ws.Range("D35") '2.883'
Set-TextValue $ws.Range("E35") '  -0.45%  '
Set-TextValue $ws.Range("D36") '0.9964'
Set-TextValue $ws.Range("E36") '  -0.37%  '
Set-TextValue $ws.Range("D37") '1.121'
Set-TextValue $ws.Range("E37") '  +0.85%  '
Set-TextValue $ws.Range("D38") '0.01950'
Set-TextValue $ws.Range("E38") '  -2.15%  '
Set-TextValue $ws.Range("D39") '2.458'
Set-TextValue $ws.Range("E39") '  +7.07%  '
Set-TextValue $ws.Range("D40") '0.05228'
Set-TextValue $ws.Range("E40") '  -1.86%  '
Set-TextValue $ws.Range("D41") '2.906'
Set-TextValue $ws.Range("E41") '  +1.58%  '
Set-TextValue $ws.Range("D42") '7.176'
Set-TextValue $ws.Range("E42") '  -2.95%  '
Set-TextValue $ws.Range("D43") '0.5254'
Set-TextValue $ws.Range("E43") '  -1.24%  '
Set-TextValue $ws.Range("D44") '0.1658'
Set-TextValue $ws.Range("E44") '  -4.31%  '
Set-TextValue $ws.Range("D45") '8.562'
Set-TextValue $ws.Range("E45") '  -2.20%  '
Set-TextValue $ws.Range("D46") '0.5048'
Set-TextValue $ws.Range("E46") '  -1.40%  '
Set-TextValue $ws.Range("D47") '10.29'
Set-TextValue $ws.Range("E47") '  -4.42%  '
Set-TextValue $ws.Range("D48") '104.28'
Set-TextValue $ws.Range("E48") '  -1.79%  '
Set-TextValue $ws.Range("D49") '0.9959'
Set-TextValue $ws.Range("E49") '  -0.42%  '
Set-TextValue $ws.Range("D50") '1.663'
Set-TextValue $ws.Range("E50") '  -2.40%  '
Set-TextValue $ws.Range("D51") '0.06311'
Set-TextValue $ws.Range("E51") '  -0.97%  '
